$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pad header row labels with leading/trailing spaces (ignore missing info on import) ---
$ws.Range("A1").Value = "Nom "
$ws.Range("B1").Value = " Prénom(s)"
$ws.Range("C1").Value = "Email "
$ws.Range("D1").Value = "Téléphone "
$ws.Range("E1").Value = "Siteweb "
$ws.Range("F1").Value = "Adresse "
$ws.Range("G1").Value = "Code Postal "
$ws.Range("H1").Value = "Ville "
$ws.Range("I1").Value = "Pays "
$ws.Range("J1").Value = "Commentaires "

# --- Add a trailing blank row (single-space placeholder cell) ---
$ws.Range("A9").Value = " "

# --- Turn the header+data range into a native Excel Table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:J9"), $null, 1)
$tbl.Name = "Table3"
$tbl.TableStyle = "TableStyleMedium24"

# --- Column widths (AutoFit-style best-fit widths captured in the diff) ---
$ws.Columns.Item(2).ColumnWidth = 12.28125
$ws.Columns.Item(4).ColumnWidth = 12.28125
$ws.Columns.Item(5).ColumnWidth = 10.140625
$ws.Columns.Item(6).ColumnWidth = 10.00390625
$ws.Columns.Item(7).ColumnWidth = 13.421875
$ws.Columns.Item(10).ColumnWidth = 15.8515625
